$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update timestamp message in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 02:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 586941
$ws.Range("C4").Value = 26641
$ws.Range("D4").Value = 36948
$ws.Range("E4").Value = 526353
$ws.Range("G4").Value = 1535
$ws.Range("H4").Value = 23640

# Row 10 - China
$ws.Range("B10").Value = 82249
$ws.Range("C10").Value = 89
$ws.Range("D10").Value = 77738
$ws.Range("E10").Value = 1170
$ws.Range("F10").Value = 116

# Row 84 - Crucero
$ws.Range("D84").Value = 639
$ws.Range("E84").Value = 61
$ws.Range("F84").Value = 7
$ws.Range("H84").Value = 12

# Row 153 - Bahamas
$ws.Range("B153").Value = 49
$ws.Range("C153").Value = 3
$ws.Range("E153").Value = 35
